$d = $word.ActiveDocument

# Paragraph 1: update title/date line (keep Normal style)
$d.Paragraphs(1).Range.Text = "המאמר היומי של אביב ומייק: 20.06.25" + [char]11 + "Evolving Deeper LLM Thinking"

# Paragraph 2: replace intro paragraph text
$d.Paragraphs(2).Range.Text = "המאמר מציג שיטה לשיפור ביצועי מודלי שפה בזמן אינפרנס(test-time compute). השיטה ממירה את בעיית החיפוש במרחב של פתרונות טקסטואליים לתהליך אבולוציוני מונחה־ביקורת, שנבנה כולו סביב יכולות הגנרטיבית וה""רפלקטיבי״ של ה-LLM עצמו. אין כאן fine-tuning או עדכון משקולות המודל אלא השיפור מתבצע באופק החישובי של ה-inference בלבד."

# Paragraph 3: replace text, then insert new Heading3 paragraph after it
$d.Paragraphs(3).Range.Text = "הנחת המוצא היא שבעיות רבות, כגון תכנון מסלול טיול או לוח זמנים לפגישות, לא ניתנות לפורמליזציה מלאה, אך כן ניתן לבחון עבורן את איכות הפתרון בעזרת פונקציית הערכה חיצונית. זה יוצר תרחיש שבו לא ניתן לייצר פתרונות ישירות על ידי אופטימיזציה מסורתית, אך כן ניתן לבצע חיפוש מונחה־הערכה. המאמר מבצע זאת באמצעות מנגנון גנטי שמתממש כולו בשפה טבעית."
$d.Paragraphs(3).Range.InsertParagraphAfter()
$d.Paragraphs(4).Style = "Heading 3"
$d.Paragraphs(4).Range.Text = "רכיב האלגוריתם 1: אוכלוסייה טקסטואלית"

# Paragraph 4 (orig) -> now paragraph 5: replace text, then insert new Heading4 paragraph after it
$d.Paragraphs(5).Range.Text = "כל פתרון מיוצג כטקסט כלומר תיאור מילולי של תוכנית פעולה. המרחב שבו מתבצע החיפוש אינו מרחב וקטורי ואינו בעל מבנה טופולוגי ברור. אין מרחק מוגדר בין שני פתרונות, ואין דרך לקבוע ""כיוון שיפור"". השיפור מבוצע באמצעות recombination(מושג מהמאמר) לשוני, כלומר  כתיבה מחדש של טקסט על בסיס טקסטים קודמים."
$d.Paragraphs(5).Range.InsertParagraphAfter()
$d.Paragraphs(6).Style = "Heading 4"
$d.Paragraphs(6).Range.Text = "רכיב האלגוריתם 2: מבנה אבולוציוני עם איים"

# Paragraph 5 (orig) -> now paragraph 7: replace text
$d.Paragraphs(7).Range.Text = "במקום אוכלוסייה אחת, האלגוריתם מחלק את מרחב הפתרונות למספר אוכלוסיות נפרדות הנקראים איים במאמר. כל אי עובר תהליך אבולוציוני עצמאי, אך כל כמה איטרציות מתבצעת ""הגירה"" של פתרונות מוצלחים בין האיים. כך נשמר איזון בין חיפוש מקומי (exploitation) לחיפוש גלובלי (exploration)."

# Paragraph 6 (orig) -> now paragraph 8: replace text (style remains Normal)
$d.Paragraphs(8).Range.Text = "רכיב האלגוריתם 3: בחירה מבוססת סלקציה רכה"

# Paragraph 7 (orig) -> now paragraph 9: replace text, then append 9 new paragraphs after it
$d.Paragraphs(9).Range.Text = "הבחירה של אילו פתרונות ישמשו הבסיס(הורים) לדור הבא אינה דטרמיניסטית. האלגוריתם בוחר פתרונות עם הסתברות שתלויה באיכותם, אך משמר גם סיכוי לבחירת פתרונות בינוניים, כדי למנוע התכנסות מוקדמת. זה יוצר מנגנון של סלקציה רכה שמאפשר לאוכלוסייה לשמור על גיוון מבני ורעיוני(קצת דומה לMCTS אבל בלי עצים)."

$insertionPoint = $d.Paragraphs(9)
$insertionPoint.Range.InsertParagraphAfter()
$insertionPoint = $d.Paragraphs(10)
$insertionPoint.Style = "Heading 4"
$insertionPoint.Range.Text = "רכיב האלגוריתם 4: recombination באמצעות שיח ביקורתי"

$insertionPoint.Range.InsertParagraphAfter()
$insertionPoint = $d.Paragraphs(11)
$insertionPoint.Style = "Normal"
$insertionPoint.Range.Text = "במקום לבצע recombination באמצעות תהליכים סינתטיים כמו דילוג על שורות או חיבור משפטים, האלגוריתם מייצר שיח פנימי בין שני ישויות קונספטואליות, מבקר ומחבר, אשר לומדים מהפידבק של פונקציית ההערכה. התוצאה היא טקסט חדש, שלא בהכרח בנוי כשילוב כלשהו של פתרונות קודמים, אלא כפרשנות מחודשת עליהם. תהליך זה חוזר על עצמו מספר פעמים בכל דור."

$insertionPoint.Range.InsertParagraphAfter()
$insertionPoint = $d.Paragraphs(12)
$insertionPoint.Style = "Normal"
$insertionPoint.Range.Text = "התהליך כולו מסתמך על פונקציית הערכה חיצונית שיכול להיות קוד, תוכנה או מודל נוסף שמספקת גם ציון איכות וגם פידבק טקסטואלי מפורש. חשוב להדגיש: המשוב אינו בהכרח מספרי בלבד, אלא יכול לכלול תיאור מפורט של תקלות או סטיות מהאילוצים, מה שמאפשר למודל להשתמש בו כחומר גלם לרפלקסיה."

$insertionPoint.Range.InsertParagraphAfter()
$insertionPoint = $d.Paragraphs(13)
$insertionPoint.Style = "Heading 3"
$insertionPoint.Range.Text = "יתרונות מבניים"

$insertionPoint.Range.InsertParagraphAfter()
$insertionPoint = $d.Paragraphs(14)
$insertionPoint.Style = "Normal"
$insertionPoint.Range.Text = "סקלביליות לבעיות לא מוגדרות היטב: מאחר והאלגוריתם פועל על טקסטים ולא על מבנים פורמליים, ניתן להפעילו גם כאשר אין תיאור פורמלי של הבעיה."

$insertionPoint.Range.InsertParagraphAfter()
$insertionPoint = $d.Paragraphs(15)
$insertionPoint.Style = "Normal"
$insertionPoint.Range.Text = "הפרדה בין גנרוט לאבלואציה: בניגוד לגישות המבוססות על התקדמות ליניארית כמו Chain-of-Thought או Reflexion, כאן יש חלוקה ברורה: המודל מייצר, ההערכה בוחנת, ואז מתבצע רה-קונפיגורציה של הפתרון." + [char]11 + "מניעת התכנסות מוקדמת: בזכות האיים, ההגרלות הרכות, וה-reset התקופתי, נמנעת קריסה מוקדמת לפתרונות לוקליים." + [char]11

$insertionPoint.Range.InsertParagraphAfter()
$insertionPoint = $d.Paragraphs(16)
$insertionPoint.Style = "Normal"
$insertionPoint.Range.Text = "האלגוריתם מאפשר ל-LLMs לחשוב לעומק לא דרך ניתוח סמנטי או לוגי של השפה, אלא דרך דינמיקה של תחרות, ביקורת, רפלקסיה והתמרה. זהו תהליך חישובי שמשתמש בשפה עצמה כחומר גלם לבניית פתרונות, ומוביל לשיפור איכותי של היכולות התכנוניות של המודל גם במצבים שבהם לא ניתן להגדיר מראש את מהות ""הפתרון הנכון""."

$insertionPoint.Range.InsertParagraphAfter()
$insertionPoint = $d.Paragraphs(17)
$insertionPoint.Style = "Normal"
$insertionPoint.Range.Text = "אם נביט בזה כתשתית רעיונית, המאמר מציע גישה כללית ל-meta-reasoning של מודלים: מערכת שמארגנת את החשיבה של המודל לא רק דרך פרומפט אלא דרך שילוב של רעיונות מתחרים שמתפתחים בהכוונת ביקורת. מדובר בתפיסה לא ליניארית של אינפרנס, כזו שמניחה שמחשבה טובה נולדת לא בבת אחת, אלא דרך אקספלורציה, שגיאות, ותיקון מצטבר."

$insertionPoint.Range.InsertParagraphAfter()
$insertionPoint = $d.Paragraphs(18)
$insertionPoint.Style = "Normal"
$insertionPoint.Range.Text = "https://arxiv.org/abs/2501.09891"
